$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 885 (everything from old row 885 downward shifts to 886..927)
$ws.Rows("885:885").Insert()

# Populate the new row. Column A holds a date formatted as plain text
# (e.g. "2026/02/26"), not a real Excel date value, so use a leading
# apostrophe to force text entry and avoid automatic date conversion.
$ws.Range("A885").Value = "'2026/02/26"
# Drop the "quote prefix" cell style the apostrophe entry introduces so
# the cell's style matches its plain, unstyled neighbours.
$ws.Range("A885").Style = "Normal"
$ws.Range("B885").Value = "木"
$ws.Range("C885").Value = 2
$ws.Range("D885").Value = 201
